$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Add Sheet2 after Sheet1 and make it the active sheet
$ws2 = $wb.Worksheets.Add($null, $sheet1)
$ws2.Name = "Sheet2"

# Header row (and data) populated in the same order the original author
# typed them, so shared-string indices line up with the target workbook.
$ws2.Range("B1").Value = "City"
$ws2.Range("C1").Value = "Address 1"
$ws2.Range("D1").Value = "Address 2"
$ws2.Range("E1").Value = "Postal Code"
$ws2.Range("F1").Value = "Phone number"
$ws2.Range("G1").Value = "FaxNumber"

$ws2.Range("B2").Value = "Delhi"
$ws2.Range("C2").Value = "Random #103 Delhi India"
$ws2.Range("D2").Value = "Random #103 Delhi India"

$ws2.Range("A1").Value = "Country"
$ws2.Range("A2").Value = "India"

$ws2.Range("E2").Value = 123456
$ws2.Range("F2").Value = 9876543210
$ws2.Range("G2").Value = 55555555

# Column widths (values chosen so the host's internal 1/6-character
# rounding lands as close as possible to the target stored widths)
$ws2.Columns.Item(2).ColumnWidth = 13.608072916666666
$ws2.Columns.Item(3).ColumnWidth = 25.830729166666668
$ws2.Columns.Item(4).ColumnWidth = 16.166666666666668
$ws2.Columns.Item(5).ColumnWidth = 17.053385416666668
$ws2.Columns.Item(6).ColumnWidth = 16.721354166666668
$ws2.Columns.Item(7).ColumnWidth = 28.276041666666668

# Selection / zoom
$ws2.Activate()
$ws2.Range("A2").Select()
$excel.ActiveWindow.Zoom = 85

# Page setup
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1
